# Consolidated Changes
# - Sheet2 (the "CNN" receptive-field calculator): remove padding (1 -> 0) on
#   the two 1x1 conv rows (23 & 24). The nout/nin/Activation_Size cells are
#   table-calculated-column formulas, so they ripple through automatically
#   on recalc.
# - Sheet5 (the experiment log table): fill in the previously-empty
#   "Base Skeleton", "Batch Normalization" and "Dropout" experiment rows
#   (3, 4, 5) with their Target / Parameters / Accuracy / Analysis write-up.
# - Active sheet/selection moves from Sheet2!O24 to Sheet5!B3.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet2")
$ws2 = $wb.Worksheets.Item("Sheet5")

# ---------------------------------------------------------------------
# Sheet2: drop padding on rows 23 & 24 from 1 -> 0 (dependent formula
# cells nout/nin/Activation_Size recompute automatically on recalc)
# ---------------------------------------------------------------------
$ws1.Range("F23").Value = 0
$ws1.Range("F24").Value = 0

# ---------------------------------------------------------------------
# Sheet5: populate row 3 - "Base Skeleton Model"
# (Analysis/F3 was authored before Target/B3, which is what drives the
# shared-string allocation order - F3's text becomes si=45, B3's si=46)
# ---------------------------------------------------------------------
$r3f = $ws2.Range("F3")
$r3f.WrapText = $true
$r3f.Borders.LineStyle = 1
$r3f.Value = "•  We have structured our model in a readable way`n•  The model is lighter with less number of parameters `n•  The performace is reduced compared to previous models. Since we have reduced model capacity, this is expected, the model has capability to learn. `n• Next, we will be tweaking this model further and increase the capacity to push it more towards the desired accuracy."

$r3b = $ws2.Range("B3")
$r3b.WrapText = $true
$r3b.Borders.LineStyle = 1
$r3b.Value = " • Get the basic skeleton interms of convolution and placement of transition blocks (max pooling, 1x1's)`n•  Reduce the number of parameters as low as possible`n•  Add GAP and remove the last BIG kernel."

$ws2.Range("C3").Value = 4572
$ws2.Range("D3").Value = 98.22
$ws2.Range("E3").Value = 98.43

$ws2.Rows.Item(3).RowHeight = 115.2

# ---------------------------------------------------------------------
# Sheet5: populate row 4 - "With Batch Normalization"
# ---------------------------------------------------------------------
$r4b = $ws2.Range("B4")
$r4b.Borders.LineStyle = 1
$r4b.Value = "•  Add Batch-norm to increase model efficiency."

$ws2.Range("C4").Value = 5088
$ws2.Range("D4").Value = 99.03
$ws2.Range("E4").Value = 99.04

$r4f = $ws2.Range("F4")
$r4f.WrapText = $true
$r4f.Borders.LineStyle = 1
$r4f.Value = "•  There is slight increase in the number of parameters, as batch norm stores a specific mean and std deviation for each layer `n • Model overfitting problem is rectified to an extent. But, we have not reached the target test accuracy 99.40%."

$ws2.Rows.Item(4).RowHeight = 86.4

# ---------------------------------------------------------------------
# Sheet5: populate row 5 - "With Dropout"
# ---------------------------------------------------------------------
$r5b = $ws2.Range("B5")
$r5b.WrapText = $true
$r5b.Borders.LineStyle = 1
$r5b.Value = "`n•  Add Regularization Dropout to each layer except last layer"

$ws2.Range("C5").Value = 5088
$ws2.Range("D5").Value = 97.94
$ws2.Range("E5").Value = 98.64

$r5f = $ws2.Range("F5")
$r5f.WrapText = $true
$r5f.Borders.LineStyle = 1
$r5f.Value = "•  There is no overfitting at all. With dropout training will be harder, because we are droping the pixels randomly.`n•  The performance has droppped, we can further improve it. `n•  But with the current capacity,not possible to push it further.We can possibly increase the capacity of the model by adding a layer after GAP! "

$ws2.Rows.Item(5).RowHeight = 115.2

# ---------------------------------------------------------------------
# Selection / active sheet: Sheet2!C24 selected (not active), Sheet5
# becomes the active sheet with B3 selected.
# ---------------------------------------------------------------------
$ws1.Range("C24").Select()
$ws2.Activate()
$ws2.Range("B3").Select()
